$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Gpha2"
$ws.Cells.Item(2, 3).Value = "Tshr"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.06840833333333333
$ws.Cells.Item(2, 8).Value = 0.205225
$ws.Cells.Item(2, 9).Value = 0.25410801246612
$ws.Cells.Item(2, 10).Value = 0.25410801246612
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.6327629999999999
$ws.Cells.Item(2, 14).Value = 1.898289
$ws.Cells.Item(2, 15).Value = 0.1382544270550543
$ws.Cells.Item(2, 16).Value = 0.1382544270550544
$ws.Cells.Item(2, 17).Value = 0.04328626222499999
$ws.Cells.Item(2, 18).Value = 0.3895763600249999
$ws.Cells.Item(2, 19).Value = 0.03513155767360203
$ws.Cells.Item(2, 20).Value = 0.03513155767360203

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Gpha2"
$ws.Cells.Item(3, 3).Value = "Tshr"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.06840833333333333
$ws.Cells.Item(3, 8).Value = 0.205225
$ws.Cells.Item(3, 9).Value = 0.25410801246612
$ws.Cells.Item(3, 10).Value = 0.25410801246612
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 2.180983666666667
$ws.Cells.Item(3, 14).Value = 6.542951
$ws.Cells.Item(3, 15).Value = 0.4765301499162115
$ws.Cells.Item(3, 16).Value = 0.4765301499162115
$ws.Cells.Item(3, 17).Value = 0.1491974576638889
$ws.Cells.Item(3, 18).Value = 1.342777118975
$ws.Cells.Item(3, 19).Value = 0.1210901292753907
$ws.Cells.Item(3, 20).Value = 0.1210901292753907

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Gpha2"
$ws.Cells.Item(4, 3).Value = "Tshr"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.06840833333333333
$ws.Cells.Item(4, 8).Value = 0.205225
$ws.Cells.Item(4, 9).Value = 0.25410801246612
$ws.Cells.Item(4, 10).Value = 0.25410801246612
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 1.444396333333334
$ws.Cells.Item(4, 14).Value = 4.333189000000001
$ws.Cells.Item(4, 15).Value = 0.3155908096798033
$ws.Cells.Item(4, 16).Value = 0.3155908096798033
$ws.Cells.Item(4, 17).Value = 0.09880874583611113
$ws.Cells.Item(4, 18).Value = 0.8892787125250001
$ws.Cells.Item(4, 19).Value = 0.08019415340030837
$ws.Cells.Item(4, 20).Value = 0.08019415340030835

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Gpha2"
$ws.Cells.Item(5, 3).Value = "Tshr"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.06840833333333333
$ws.Cells.Item(5, 8).Value = 0.205225
$ws.Cells.Item(5, 9).Value = 0.25410801246612
$ws.Cells.Item(5, 10).Value = 0.25410801246612
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.3186579999999999
$ws.Cells.Item(5, 14).Value = 0.9559739999999999
$ws.Cells.Item(5, 15).Value = 0.06962461334893082
$ws.Cells.Item(5, 16).Value = 0.06962461334893082
$ws.Cells.Item(5, 17).Value = 0.02179886268333333
$ws.Cells.Item(5, 18).Value = 0.19618976415
$ws.Cells.Item(5, 19).Value = 0.0176921721168189
$ws.Cells.Item(5, 20).Value = 0.01769217211681889

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Gpha2"
$ws.Cells.Item(6, 3).Value = "Tshr"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 0.2008013333333333
$ws.Cells.Item(6, 8).Value = 0.6024039999999999
$ws.Cells.Item(6, 9).Value = 0.7458919875338801
$ws.Cells.Item(6, 10).Value = 0.74589198753388
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.6327629999999999
$ws.Cells.Item(6, 14).Value = 1.898289
$ws.Cells.Item(6, 15).Value = 0.1382544270550543
$ws.Cells.Item(6, 16).Value = 0.1382544270550544
$ws.Cells.Item(6, 17).Value = 0.1270596540839999
$ws.Cells.Item(6, 18).Value = 1.143536886756
$ws.Cells.Item(6, 19).Value = 0.1031228693814523
$ws.Cells.Item(6, 20).Value = 0.1031228693814523

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Gpha2"
$ws.Cells.Item(7, 3).Value = "Tshr"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 0.6666666666666666
$ws.Cells.Item(7, 7).Value = 0.2008013333333333
$ws.Cells.Item(7, 8).Value = 0.6024039999999999
$ws.Cells.Item(7, 9).Value = 0.7458919875338801
$ws.Cells.Item(7, 10).Value = 0.74589198753388
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 2.180983666666667
$ws.Cells.Item(7, 14).Value = 6.542951
$ws.Cells.Item(7, 15).Value = 0.4765301499162115
$ws.Cells.Item(7, 16).Value = 0.4765301499162115
$ws.Cells.Item(7, 17).Value = 0.4379444282448888
$ws.Cells.Item(7, 18).Value = 3.941499854204
$ws.Cells.Item(7, 19).Value = 0.3554400206408208
$ws.Cells.Item(7, 20).Value = 0.3554400206408208

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Gpha2"
$ws.Cells.Item(8, 3).Value = "Tshr"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 0.6666666666666666
$ws.Cells.Item(8, 7).Value = 0.2008013333333333
$ws.Cells.Item(8, 8).Value = 0.6024039999999999
$ws.Cells.Item(8, 9).Value = 0.7458919875338801
$ws.Cells.Item(8, 10).Value = 0.74589198753388
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 1.444396333333334
$ws.Cells.Item(8, 14).Value = 4.333189000000001
$ws.Cells.Item(8, 15).Value = 0.3155908096798033
$ws.Cells.Item(8, 16).Value = 0.3155908096798033
$ws.Cells.Item(8, 17).Value = 0.2900367095951111
$ws.Cells.Item(8, 18).Value = 2.610330386356
$ws.Cells.Item(8, 19).Value = 0.235396656279495
$ws.Cells.Item(8, 20).Value = 0.2353966562794949

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Gpha2"
$ws.Cells.Item(9, 3).Value = "Tshr"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 0.6666666666666666
$ws.Cells.Item(9, 7).Value = 0.2008013333333333
$ws.Cells.Item(9, 8).Value = 0.6024039999999999
$ws.Cells.Item(9, 9).Value = 0.7458919875338801
$ws.Cells.Item(9, 10).Value = 0.74589198753388
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.3186579999999999
$ws.Cells.Item(9, 14).Value = 0.9559739999999999
$ws.Cells.Item(9, 15).Value = 0.06962461334893082
$ws.Cells.Item(9, 16).Value = 0.06962461334893082
$ws.Cells.Item(9, 17).Value = 0.06398695127733331
$ws.Cells.Item(9, 18).Value = 0.5758825614959999
$ws.Cells.Item(9, 19).Value = 0.05193244123211192
$ws.Cells.Item(9, 20).Value = 0.05193244123211192
